$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '95.886.66'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '3.621.10'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '2.74'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +29.09%  '
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '224.22'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -4.68%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '640.84'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.06%  '
$ws.Range('E8').Value = '  -2.51%  '
$ws.Range('E9').Value = '  +11.06%  '
$ws.Range('E10').Value = '  +0.07%  '
$ws.Range('D11').Value = '3.617.57'
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '48.37'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +9.08%  '
$ws.Range('E13').Value = '  +4.59%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000291'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -5.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.55'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.63%  '
$ws.Range('E16').Value = '  -1.70%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '24.75'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +33.61%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '95.645.87'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.01'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.90'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +7.35%  '
$ws.Range('D21').Value = '3.624.82'
$ws.Range('E21').Value = '  -1.66%  '
$ws.Range('E22').Value = '  +47.18%  '
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '137.44'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +24.86%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '524.96'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.72%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.27'
$ws.Range('D26').Style = "Normal"
$ws.Range('E27').Value = '  -7.60%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.87'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.15%  '
$ws.Range('D29').Value = '3.802.88'
$ws.Range('E29').Value = '  -2.13%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '12.95'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '13.26'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +6.05%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.13'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.71%  '
$ws.Range('E33').Value = '  +0.14%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.636'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +8.08%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '33.40'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.77%  '
$ws.Range('E36').Value = '  -1.76%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.83'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.538'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +9.83%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '7.26'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.15%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '591.97'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -5.95%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0535'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +20.10%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.36'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.76%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '41.38'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.65%  '
$ws.Range('E46').Value = '  +5.72%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.159'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.25%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.99'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.26'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +6.80%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '237.62'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +16.63%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.32'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.98%  '
